$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Add the new hidden "DropdownOptions" sheet right after Sheet1.
$dropdownSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$dropdownSheet.Name = "DropdownOptions"

# 2. Populate the dropdown option values (A1:A7).
$options = @("0% - 10%", "11% - 25%", "26% - 50%", "51% - 75%", "76% - 90%", "91% - 99%", "100%")
for ($i = 0; $i -lt $options.Length; $i++) {
    $cell = $dropdownSheet.Cells.Item($i + 1, 1)
    $cell.NumberFormat = "@"
    $cell.Value = $options[$i]
    $cell.Style = "Normal"
}

# 3. Hide the helper sheet.
$dropdownSheet.Visible = $false

# 4. Add the new "Status as of July 4, 2025" header column on Sheet1.
$ws.Range("AA1").Value = "Status as of July 4, 2025"

# 5. Clear out stray empty typed cells that Excel collapses away on save.
$ws.Range("Q5").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("Z8").ClearContents()
$ws.Range("Z14").ClearContents()
$ws.Range("Z15").ClearContents()

# 6. Apply list data validation (dropdown) to AA2:AA17 sourced from DropdownOptions.
$rng = $ws.Range("AA2:AA17")
$rng.Validation.Add(3, 1, 1, "=DropdownOptions!`$A`$1:`$A`$7")
$rng.Validation.IgnoreBlank = $true
$rng.Validation.InCellDropdown = $true
$rng.Validation.ShowInput = $false
$rng.Validation.ShowError = $false

# 7. Re-activate Sheet1 so it remains the active sheet/tab (matches the original file).
$ws.Activate()
